$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("NCTId"), shifting old B..I to C..J.
$ws.Columns("B:B").Insert()

# New header for the inserted column: string version of "statut".
$ws.Cells.Item(1, 2).Value = "status_label"

# Target data (row, status_label, NCTId, completion_year, clinical_trial_title, acronym)
# Rows were also re-ordered during the edit, so we rewrite B2:G20 fully with final values.
$rows = @(
    @(2, "rouge", "NCT02016222", "2016", "Tear Analysis in the Diagnosis of Primary Progressive Forms of Multiple Sclerosis", "LARMES"),
    @(3, "rouge", "NCT02276924", "2016", "Diagnostic Relevance of Laser Confocal Microscopy During Reno-ureteroscopy in the Context of the Screening and Follow-up of Upper Urinary Tract Tumors", "UROVISIO"),
    @(4, "rouge", "NCT02765529", "2016", "Research of Biomarkers of Air Pollutants Exposure", "AEROTOX-2"),
    @(5, "rouge", "NCT02005679", "2017", "Transposition of the Mini-Mental State in Sign Language", "MMS-LS"),
    @(6, "rouge", "NCT02043964", "2018", "Tear Analysis by Isoelectric Focusing in Clinically Isolated Syndrome as Multiple Sclerosis Criterion Among Patients With or Without Lesions at the Magnetic Resonance Imaging (Oligoclonal Profile of Tears)", "POLAR"),
    @(7, "rouge", "NCT02297165", "2018", "Efficacy Study of an Olfactory Stimulation Program in Relaxing Environment for the Recuperation of Autobiographical Memories in Anorexic Patients", "OLFANOR"),
    @(8, "rouge", "NCT02669160", "2018", "Tolerance of a Motorized Orthosis Reproducing Walking Movement Versus Conventional Passive Standing-up Devices in Children With Cerebral Palsy : A Non-inferiority, Randomised, Multicenter, Controlled Trial", "EOMEC/CP"),
    @(9, "rouge", "NCT03393299", "2018", "Impact of the Systematic Use of the Criteria STOPP/START in Short Stay Geriatric: Study of Superiority, Randomized, Controlled, Prospective, Single Blind", "REVOR"),
    @(10, "rouge", "NCT03666793", "2019", "Comprehensive Management of Drug Prescriptions Throughout the Elderly Person's Hospital Care, From Hospital to Home: Impact on Readmission at 30 Days After Delivery", "OPTISORT"),
    @(11, "rouge", "NCT02733900", "2020", "Study of Modifications of the Composition and Structure in the Aseptic Osteonecrosis of the Femoral Head and Etiopathogenic MRI Correlations", "ONTF"),
    @(12, "rouge", "NCT03727217", "2020", "Performance of Ultrasound in the Early Diagnosis of Vocal Cords Paralysis After Thyroidectomy or Parathyroidectomy (PECV)", "PECV"),
    @(13, "rouge", "NCT04224077", "2020", "Optimisation of Diffusion Tensor Sequences (DTI) for Study of Lumbar Roots by Magnetic Resonance Imaging (MRI) : a Feasibility Study", "OPTI-DTI"),
    @(14, "rouge", "NCT03976674", "2020", "Evaluation of a Preoperative Cognitive Behavioural Therapy (CBT) Program Based on Self-determination Theory for Bariatric Surgery Candidates : an Open-label Controlled, Randomized, Superiority Study", "ACRoBAT"),
    @(15, "rouge", "NCT02862379", "2021", "Evaluation of a Personalized Rehabilitation Program for Elderly Patients That Fall : Impact on the Fear of Falling", "CHUTE"),
    @(16, "rouge", "NCT03162341", "2021", "Study of the Correlation Between UltraSonography and Dual-Energy Computed Tomography Assessment of Urate Deposit in Urate Lowering Therapy Initiators", "GOUT"),
    @(17, "rouge", "NCT04699136", "2022", "Use and Validation of the 6-minute Stepper Test in in Patients With Cardiac Pathologies Needing Rehabilitation and Reeducation (CVRR)", "STEPPER"),
    @(18, "rouge", "NCT04169477", "2022", "Superiority, Prospective, Multicentric, Randomized, Single-blind, Cross-over Study Comparing 2 Modes of Transcutaneous Electrical Nerve Stimulation (TENS) in Chronic Neuropathic Radiculalgia", "CROSS-TENS"),
    @(19, "rouge", "NCT03128905", "2022", "Colchicine or Prednisone for the Treatment of Acute Calcium Pyrophosphate Deposition (CPPD) Arthritis: Open-label, Randomized, Multicenter, Equivalence Trial of Efficacy and Safety", "COLCHICORT"),
    @(20, "rouge", "NCT02819037", "2022", "Modification of Digestive Flora After Gastric Bypass : Pilot Study on Microbial Overgrowth Using Gas Chromatography", "SIBOB")
)

foreach ($row in $rows) {
    $r = $row[0]
    $statusLabel = $row[1]
    $nctId = $row[2]
    $year = $row[3]
    $title = $row[4]
    $acronym = $row[5]

    $ws.Cells.Item($r, 2).Value = $statusLabel
    $ws.Cells.Item($r, 3).Value = $nctId
    # eudraCT (column D) stays empty for every row.
    $ws.Cells.Item($r, 4).Value = ""
    # completion_year looks numeric; force text storage with a leading apostrophe,
    # matching the original inlineStr/text representation ("2016", not 2016).
    $ws.Cells.Item($r, 5).Value = "'" + $year
    $ws.Cells.Item($r, 6).Value = $title
    $ws.Cells.Item($r, 7).Value = $acronym
    # results_1y / results_3y / results booleans are unchanged (all FALSE),
    # they now simply live in columns H / I / J after the column insert.
    $ws.Cells.Item($r, 8).Value = $false
    $ws.Cells.Item($r, 9).Value = $false
    $ws.Cells.Item($r, 10).Value = $false
}
